$wb = $excel.ActiveWorkbook

# ---- Sheet: Summary ----
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.402135231316726
$ws1.Range("C2").Value = 0.07458563535911603
$ws1.Range("D2").Value = 0.9642857142857143
$ws1.Range("E2").Value = 0.1384615384615385
$ws1.Range("F2").Value = 0.2848101265822785
$ws1.Range("G2").Value = 0.6610169491525424
$ws1.Range("H2").Value = 0.8012974852862493
$ws1.Range("I2").Value = 27
$ws1.Range("J2").Value = 335
$ws1.Range("K2").Value = 199
$ws1.Range("L2").Value = 1

# ---- Sheet: Classification Report ----
$ws2 = $wb.Worksheets.Item("Classification Report")

$ws2.Range("B2").Value = 0.995
$ws2.Range("C2").Value = 0.3726591760299626
$ws2.Range("D2").Value = 0.5422343324250681

$ws2.Range("B3").Value = 0.07458563535911603
$ws2.Range("C3").Value = 0.9642857142857143
$ws2.Range("D3").Value = 0.1384615384615385

$ws2.Range("B4").Value = 0.402135231316726
$ws2.Range("C4").Value = 0.402135231316726
$ws2.Range("D4").Value = 0.402135231316726
$ws2.Range("E4").Value = 0.402135231316726

$ws2.Range("B5").Value = 0.534792817679558
$ws2.Range("C5").Value = 0.6684724451578384
$ws2.Range("D5").Value = 0.3403479354433033

$ws2.Range("B6").Value = 0.9491430565659348
$ws2.Range("C6").Value = 0.402135231316726
$ws2.Range("D6").Value = 0.5221175384197676

# ---- Sheet: Confusion Matrix ----
$ws3 = $wb.Worksheets.Item("Confusion Matrix")

$ws3.Range("B2").Value = 199
$ws3.Range("C2").Value = 335

$ws3.Range("B3").Value = 1
$ws3.Range("C3").Value = 27
